$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that swap values between row 2 and row 3: D, L, M, N, O, P, S
$cols = @("D", "L", "M", "N", "O", "P", "S")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $v2 = $cell2.Value2
    $v3 = $cell3.Value2
    $cell2.Value2 = $v3
    $cell3.Value2 = $v2
}
